$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 6).Value = 3.2
$ws.Cells.Item(2, 7).Value = 3.25
$ws.Cells.Item(2, 8).Value = 2.44
$ws.Cells.Item(2, 15).Value = 1.33
$ws.Cells.Item(2, 16).Value = 1.9
$ws.Cells.Item(2, 17).Value = 2.02
$ws.Cells.Item(2, 21).Value = 2.18
$ws.Cells.Item(2, 22).Value = 1.68
$ws.Cells.Item(2, 23).Value = 1.44
$ws.Cells.Item(2, 40).Value = 34
$ws.Cells.Item(4, 12).Value = 1.35
$ws.Cells.Item(4, 25).Value = 21
$ws.Cells.Item(4, 29).Value = 10.5
$ws.Cells.Item(4, 30).Value = 21
$ws.Cells.Item(4, 31).Value = 60
$ws.Cells.Item(4, 33).Value = 11
$ws.Cells.Item(4, 34).Value = 21
$ws.Cells.Item(4, 36).Value = 25
$ws.Cells.Item(4, 37).Value = 23
$ws.Cells.Item(4, 38).Value = 38
$ws.Cells.Item(4, 40).Value = 12
$ws.Cells.Item(5, 6).Value = 1.69
$ws.Cells.Item(5, 11).Value = 3.7
$ws.Cells.Item(5, 12).Value = 1.54
$ws.Cells.Item(5, 16).Value = 1.52
$ws.Cells.Item(5, 29).Value = 10.5
$ws.Cells.Item(6, 10).Value = 2.72
$ws.Cells.Item(6, 11).Value = 3.15
$ws.Cells.Item(6, 22).Value = 1.34
$ws.Cells.Item(6, 23).Value = 1.55
$ws.Cells.Item(7, 8).Value = 2.08
$ws.Cells.Item(7, 10).Value = 2.66
$ws.Cells.Item(7, 12).Value = 1.68
$ws.Cells.Item(7, 29).Value = 8.6
$ws.Cells.Item(8, 7).Value = 2.12
$ws.Cells.Item(8, 8).Value = 4.1
$ws.Cells.Item(8, 9).Value = 4.7
$ws.Cells.Item(8, 12).Value = 1.42
$ws.Cells.Item(8, 13).Value = 1.08
$ws.Cells.Item(8, 18).Value = 1.29
$ws.Cells.Item(8, 19).Value = 3.65
$ws.Cells.Item(8, 20).Value = 1.84
$ws.Cells.Item(8, 21).Value = 1.96
$ws.Cells.Item(8, 22).Value = 1.27
$ws.Cells.Item(8, 23).Value = 1.9
$ws.Cells.Item(8, 24).Value = 1000
$ws.Cells.Item(8, 25).Value = 1000
$ws.Cells.Item(8, 26).Value = 980
$ws.Cells.Item(8, 28).Value = 1000
$ws.Cells.Item(8, 29).Value = 8.4
$ws.Cells.Item(8, 30).Value = 1000
$ws.Cells.Item(8, 31).Value = 60
$ws.Cells.Item(8, 32).Value = 1000
$ws.Cells.Item(8, 33).Value = 1000
$ws.Cells.Item(8, 34).Value = 1000
$ws.Cells.Item(8, 35).Value = 70
$ws.Cells.Item(8, 36).Value = 980
$ws.Cells.Item(8, 37).Value = 980
$ws.Cells.Item(8, 39).Value = 140
$ws.Cells.Item(8, 40).Value = 1000
$ws.Cells.Item(9, 7).Value = 2.18
$ws.Cells.Item(9, 8).Value = 4.6
$ws.Cells.Item(9, 10).Value = 2.78
$ws.Cells.Item(9, 11).Value = 3.35
$ws.Cells.Item(9, 12).Value = 1.56
$ws.Cells.Item(9, 14).Value = 2.2
$ws.Cells.Item(9, 15).Value = 1.67
$ws.Cells.Item(9, 17).Value = 2.76
$ws.Cells.Item(9, 20).Value = 2.46
$ws.Cells.Item(9, 21).Value = 1.56
$ws.Cells.Item(9, 23).Value = 1.86
$ws.Cells.Item(9, 38).Value = 1000
$ws.Cells.Item(10, 12).Value = 1.35
$ws.Cells.Item(10, 14).Value = 3.95
$ws.Cells.Item(10, 18).Value = 1.41
$ws.Cells.Item(10, 19).Value = 2.98
$ws.Cells.Item(10, 20).Value = 1.79
$ws.Cells.Item(10, 21).Value = 2.02
$ws.Cells.Item(11, 6).Value = 1.8
$ws.Cells.Item(11, 9).Value = 5.3
$ws.Cells.Item(11, 12).Value = 1.36
$ws.Cells.Item(11, 13).Value = 1.05
$ws.Cells.Item(11, 14).Value = 3.8
$ws.Cells.Item(11, 16).Value = 1.98
$ws.Cells.Item(11, 18).Value = 1.38
$ws.Cells.Item(11, 19).Value = 3.05
$ws.Cells.Item(11, 20).Value = 1.74
$ws.Cells.Item(11, 21).Value = 2.06
$ws.Cells.Item(11, 22).Value = 1.21
$ws.Cells.Item(11, 23).Value = 1.97
$ws.Cells.Item(11, 24).Value = 22
$ws.Cells.Item(11, 28).Value = 11.5
$ws.Cells.Item(11, 29).Value = 11
$ws.Cells.Item(11, 30).Value = 27
$ws.Cells.Item(11, 31).Value = 75
$ws.Cells.Item(11, 32).Value = 14.5
$ws.Cells.Item(11, 33).Value = 1000
$ws.Cells.Item(11, 35).Value = 75
$ws.Cells.Item(11, 36).Value = 25
$ws.Cells.Item(11, 40).Value = 14.5
$ws.Cells.Item(12, 12).Value = 1.25
$ws.Cells.Item(12, 18).Value = 1.57
$ws.Cells.Item(12, 19).Value = 2.4
$ws.Cells.Item(12, 20).Value = 1.79
$ws.Cells.Item(12, 28).Value = 32
$ws.Cells.Item(19, 12).Value = 1.32
$ws.Cells.Item(19, 18).Value = 1.54
$ws.Cells.Item(20, 19).Value = 3.8
$ws.Cells.Item(20, 30).Value = 21
$ws.Cells.Item(21, 9).Value = 5.1
$ws.Cells.Item(21, 15).Value = 1.01
$ws.Cells.Item(22, 12).Value = 1.53
$ws.Cells.Item(22, 13).Value = 1.11
$ws.Cells.Item(22, 14).Value = 2.8
$ws.Cells.Item(22, 18).Value = 1.22
$ws.Cells.Item(22, 19).Value = 4.5
$ws.Cells.Item(22, 20).Value = 1.96
$ws.Cells.Item(22, 21).Value = 1.87
$ws.Cells.Item(22, 22).Value = 1.41
$ws.Cells.Item(22, 24).Value = 10.5
$ws.Cells.Item(22, 25).Value = 10.5
$ws.Cells.Item(22, 26).Value = 22
$ws.Cells.Item(22, 27).Value = 70
$ws.Cells.Item(22, 28).Value = 9.199999999999999
$ws.Cells.Item(22, 29).Value = 7.6
$ws.Cells.Item(22, 30).Value = 15
$ws.Cells.Item(22, 31).Value = 980
$ws.Cells.Item(22, 32).Value = 17.5
$ws.Cells.Item(22, 33).Value = 13.5
$ws.Cells.Item(22, 34).Value = 22
$ws.Cells.Item(22, 36).Value = 980
$ws.Cells.Item(22, 37).Value = 980
$ws.Cells.Item(22, 38).Value = 60
$ws.Cells.Item(22, 39).Value = 180
$ws.Cells.Item(22, 40).Value = 48
$ws.Cells.Item(22, 41).Value = 65
$ws.Cells.Item(23, 7).Value = 3.85
$ws.Cells.Item(23, 9).Value = 2.4
$ws.Cells.Item(23, 12).Value = 1.39
$ws.Cells.Item(23, 13).Value = 1.06
$ws.Cells.Item(23, 17).Value = 1.81
$ws.Cells.Item(23, 22).Value = 1.71
$ws.Cells.Item(24, 9).Value = 2.58
$ws.Cells.Item(24, 10).Value = 2.98
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 25).Value = 6.6
